$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.189845085144043
$ws.Range("B1").Value = 2.575979471206665
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.203083276748657
$ws.Range("E1").Value = 1.179607152938843
